# "Drop in files from RMI script"
# - Remove the obsolete "Texas Notes" worksheet entirely.
# - Update the DR discount rate from 5.87% back down to 3%.
# - Refresh the on-screen selections left behind in the "About" and "DR" sheets.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Drop the "Texas Notes" sheet (sheetId 3) - no longer needed.
$wb.Worksheets("Texas Notes").Delete()

# Update the discount rate value on the DR sheet.
$wsDR = $wb.Worksheets("DR")
$wsDR.Activate()
$wsDR.Range("B2").Value = 0.03
$wsDR.Range("B1").Select()

# Restore the "About" sheet as the active tab with its refreshed selection.
$wsAbout = $wb.Worksheets("About")
$wsAbout.Activate()
$wsAbout.Range("A16:A18").Select()
